# Add a "Run 50" results column. The previously-last data column (AZ),
# which held the "Mean" header/values, is repurposed to hold the new
# run's values, and a brand-new column (BA) is appended to hold the
# "Mean" header and the recomputed mean (now averaging 51 runs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value used for the "Run 50" column on every data row (2-14), and
# the resulting new mean (now averaging 51 runs instead of 50).
$newRunValue = 111430409228.6299
$newMeanValue = 94668309596.25479

# Header row: AZ1 currently holds "Mean" -> becomes "Run 50"; the new
# BA1 becomes "Mean", formatted like the other header cells (copy the
# bold/centered/bordered style from AZ1).
$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)

# Data rows 2-14: AZ column becomes the new Run 50 value, BA column
# becomes the recomputed mean.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = $newRunValue
    $ws.Cells.Item($r, 53).Value = $newMeanValue
}
